# Apply the "ActorTable" updates:
#  - BigBatSuccubus attackDelay tuned 0.7 -> 0.77
#  - new "flying|Bool" column (O) added, defaulted to FALSE for all existing rows
#  - new JellyFishGirl row (row 6) fully populated

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ActorTable")

# New column header (this is the first brand-new shared string -> keeps
# sharedStrings.xml ordering identical to the authored edit)
$ws.Range("O1").Value = "flying|Bool"

# BigBat Succubus (row 4) attack-speed tweak
$ws.Range("G4").Value = 0.77

# row 2's multiAtk value was also retuned alongside the column addition
$ws.Range("F2").Value = 0.95

# JellyFishGirl registration (row 6)
$ws.Range("B6").Value = "CharName_JellyFishGirl"
$ws.Range("C6").Value = "CharDesc_JellyFishGirl"
$ws.Range("E6").Value = 0.87
$ws.Range("F6").Value = 0.82
$ws.Range("G6").Value = 0.83
$ws.Range("I6").Value = 1
$ws.Range("L6").Value = "JellyFishGirl"
$ws.Range("M6").Value = "Portrait_JellyFishGirl"
$ws.Range("N6").Value = 0.05

# Populate the new flying|Bool column for every data row (all default to FALSE)
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 15).Value = $false
}

Write-Host "ActorTable updated"
